$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("O1").Value = "Serviced by "

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = ""
}
